$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 'As I leave her house and turn back to look at the building, I feel a sense of desolation.'
$ws.Range("G3").Value = '…Mana Mizuki…chan….'
$ws.Range("G4").Value = 'I recall the memory of that childish girl.'
$ws.Range("G5").Value = 'I remember how she had smiled kindly at me at the station that one time, but something felt different…'
$ws.Range("G6").Value = 'Like she was a different person…'
$ws.Range("G7").Value = 'I don''t think that''s really the case, but still….'
$ws.Range("G8").Value = '…I''m home….'
$ws.Range("G9").Value = 'Oh, there''s a message on the answering machine.'
$ws.Range("G10").Value = 'Beep….'
$ws.Range("D11").Value = 'Answering Machine'
$ws.Range("D12").Value = 'Answering Machine'
$ws.Range("D15").Value = 'Answering Machine'
$ws.Range("D17").Value = 'Answering Machine'
$ws.Range("D21").Value = 'Answering Machine'
$ws.Range("G11").Value = 'This is the Home Tutor Center.'
$ws.Range("G12").Value = 'Regarding Mizuki-san''s request, we appreciate your acceptance, and look forward to working with you.'
$ws.Range("G13").Value = 'Eh…...?'
$ws.Range("G14").Value = 'Did I make such a call?'
$ws.Range("G15").Value = 'Well then, we look forward to working with you until March of next year.'
$ws.Range("G16").Value = '…Ah, wait, could it be that Mana-chan did this on her own…?'
$ws.Range("G17").Value = '…Also, this is a request from her parents, but if Mizuki-san does not show enough enthusiasm in her studies, we''d like you to report it to us.'
$ws.Range("G18").Value = '…I see.'
$ws.Range("G19").Value = 'I am a watchdog, huh….'
$ws.Range("G20").Value = 'It looks like her story wasn''t just a case of excessive self-consciousness.'
$ws.Range("G21").Value = 'Well then, thank you for your time…'
$ws.Range("G22").Value = 'Beep….'
$ws.Range("G23").Value = 'I don''t know what''s going on, but either way, Mana-chan doesn''t seem to be in a very good situation.'
$ws.Range("G24").Value = 'Even if I refuse now, another "watchdog" will just show up eventually.'
$ws.Range("G25").Value = 'I was told it''s okay for me to come as I please, so maybe I don''t need to refuse after all….'

$ws.Range("G10").Select()
